$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New "pre" / "post" header labels in column A / G ---
$ws.Range("A1").Value = "pre"
$ws.Range("G1").Value = "post"

# --- 2. Duplicate the B:E pinout table into H:K ("post" table) ---
# Rows that carry the full Function/Timer/Channel/Pin record (B:E -> H:K)
$fullRows = @(1,2,3,4,6,7,9,11,12,13,14,15,16,18,20,21,22,23,24,25,26)
# Blank separator rows that only carry styled Timer/Channel cells (C:D -> I:J)
$sepRows  = @(5,8,10,17,19)

foreach ($r in $fullRows) {
    $ws.Range("B$r`:E$r").Copy() | Out-Null
    $ws.Range("H$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("B$r`:E$r").Copy() | Out-Null
    $ws.Range("H$r").PasteSpecial(-4163) | Out-Null   # xlPasteValues
}
foreach ($r in $sepRows) {
    $ws.Range("C$r`:D$r").Copy() | Out-Null
    $ws.Range("I$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = 0

# --- 3. Fix the ENC2x / ENC3x pin-channel mix-up in the original (B:E) table ---
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "PA8"

$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "PA9"

$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = "PA4"

$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = "PA6"

# --- swap M1D / M3D pins ---
$ws.Range("E20").Value = "PA7"
$ws.Range("E22").Value = "PA1"

# --- 4. Apply the same corrected values to the new "post" table (H:K) ---
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = "PA8"

$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = "PA9"

$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = "PA4"

$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = "PA6"

$ws.Range("K20").Value = "PA7"
$ws.Range("K22").Value = "PA1"

# --- 5. "post" table traces for S1PWM / S2PWM moved to timer 15 ---
$ws.Range("I6").Value = 15
$ws.Range("J6").Value = 1
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = 2

# --- 6. Selection / view bookkeeping ---
$ws.Range("N13").Select() | Out-Null
